# Update the cached "datetimeFigureOut" date field text shown in the
# Date Placeholder of the slide master and every slide layout from
# 11/22/2022 to 12/12/2022.

$p = $ppt.ActivePresentation
$oldDate = "11/22/2022"
$newDate = "12/12/2022"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        if ($shp.Name -like "Date Placeholder*") {
            $isDatePlaceholder = $true
        } elseif ($shp.Type -eq 14) {
            try {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePlaceholder = $true
                }
            } catch {
            }
        }
        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Custom (slide) Layout belonging to the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    Update-DatePlaceholder $lay.Shapes
}
